$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.2103559870550162
$ws.Range("C2").Value = 0.5210355987055016
$ws.Range("J2").Value = 0.01941747572815534
$ws.Range("P2").Value = 0.1553398058252427
$ws.Range("S2").Value = 0.09385113268608414
# Row 3
$ws.Range("B3").Value = 0.005952380952380952
$ws.Range("C3").Value = 0.02380952380952381
$ws.Range("J3").Value = 0.02380952380952381
$ws.Range("P3").Value = 0.6726190476190477
$ws.Range("S3").Value = 0.2738095238095238
# Row 4
$ws.Range("J4").Value = 0.09302325581395349
$ws.Range("P4").Value = 0.6744186046511628
$ws.Range("S4").Value = 0.2325581395348837
# Row 6
$ws.Range("B6").Value = 0.06222222222222222
$ws.Range("D6").Value = 0.01333333333333333
$ws.Range("F6").Value = 0.05333333333333334
$ws.Range("J6").Value = 0.2488888888888889
$ws.Range("O6").Value = 0.03555555555555556
$ws.Range("Q6").Value = 0.1644444444444444
$ws.Range("R6").Value = 0.05777777777777778
$ws.Range("S6").Value = 0.3644444444444445
# Row 7
$ws.Range("B7").Value = 0.08163265306122448
$ws.Range("D7").Value = 0.01530612244897959
$ws.Range("F7").Value = 0.07142857142857142
$ws.Range("J7").Value = 0.1071428571428571
$ws.Range("O7").Value = 0.01530612244897959
$ws.Range("Q7").Value = 0.1530612244897959
$ws.Range("R7").Value = 0.08163265306122448
$ws.Range("S7").Value = 0.4744897959183674
# Row 8
$ws.Range("B8").Value = 0.1060240963855422
$ws.Range("D8").Value = 0.009638554216867471
$ws.Range("E8").Value = 0.002409638554216868
$ws.Range("F8").Value = 0.05542168674698795
$ws.Range("J8").Value = 0.0963855421686747
$ws.Range("O8").Value = 0.03373493975903614
$ws.Range("Q8").Value = 0.1421686746987952
$ws.Range("R8").Value = 0.1373493975903614
$ws.Range("S8").Value = 0.4168674698795181
# Row 9
$ws.Range("B9").Value = 0.1005025125628141
$ws.Range("D9").Value = 0.01507537688442211
$ws.Range("F9").Value = 0.06030150753768844
$ws.Range("J9").Value = 0.1306532663316583
$ws.Range("O9").Value = 0.01507537688442211
$ws.Range("Q9").Value = 0.1708542713567839
$ws.Range("R9").Value = 0.05527638190954774
$ws.Range("S9").Value = 0.4522613065326633
# Row 10
$ws.Range("B10").Value = 0.1125
$ws.Range("D10").Value = 0.0234375
$ws.Range("E10").Value = 0.00390625
$ws.Range("F10").Value = 0.06953125
$ws.Range("J10").Value = 0.10234375
$ws.Range("O10").Value = 0.021875
$ws.Range("Q10").Value = 0.21015625
$ws.Range("R10").Value = 0.07109375
$ws.Range("S10").Value = 0.38515625
# Row 11
$ws.Range("F11").Value = 0.003333333333333334
$ws.Range("G11").Value = 0.1533333333333333
$ws.Range("J11").Value = 0.08666666666666667
$ws.Range("K11").Value = 0.1966666666666667
$ws.Range("L11").Value = 0.5466666666666666
$ws.Range("S11").Value = 0.01333333333333333
# Row 12
$ws.Range("G12").Value = 0.7470588235294118
$ws.Range("J12").Value = 0.2
$ws.Range("K12").Value = 0.005882352941176471
$ws.Range("L12").Value = 0.01764705882352941
$ws.Range("S12").Value = 0.02941176470588235
# Row 13
$ws.Range("G13").Value = 0.6904761904761905
$ws.Range("J13").Value = 0.2857142857142857
$ws.Range("S13").Value = 0.02380952380952381
# Row 15
$ws.Range("F15").Value = 0.01869158878504673
$ws.Range("H15").Value = 0.08411214953271028
$ws.Range("I15").Value = 0.03271028037383177
$ws.Range("J15").Value = 0.3925233644859813
$ws.Range("K15").Value = 0.06542056074766354
$ws.Range("M15").Value = 0.01401869158878505
$ws.Range("O15").Value = 0.07943925233644859
$ws.Range("S15").Value = 0.3130841121495327
# Row 16
$ws.Range("F16").Value = 0.02645502645502645
$ws.Range("H16").Value = 0.1798941798941799
$ws.Range("I16").Value = 0.1111111111111111
$ws.Range("J16").Value = 0.4074074074074074
$ws.Range("K16").Value = 0.1005291005291005
$ws.Range("M16").Value = 0.02116402116402116
$ws.Range("O16").Value = 0.02645502645502645
$ws.Range("S16").Value = 0.126984126984127
# Row 17
$ws.Range("F17").Value = 0.00936768149882904
$ws.Range("H17").Value = 0.1615925058548009
$ws.Range("I17").Value = 0.08196721311475409
$ws.Range("J17").Value = 0.4426229508196721
$ws.Range("K17").Value = 0.1007025761124122
$ws.Range("M17").Value = 0.0117096018735363
$ws.Range("O17").Value = 0.05386416861826698
$ws.Range("S17").Value = 0.1381733021077283
# Row 18
$ws.Range("F18").Value = 0.03191489361702127
$ws.Range("H18").Value = 0.0851063829787234
$ws.Range("I18").Value = 0.09574468085106383
$ws.Range("J18").Value = 0.425531914893617
$ws.Range("K18").Value = 0.1595744680851064
$ws.Range("M18").Value = 0.01063829787234043
$ws.Range("O18").Value = 0.05319148936170213
$ws.Range("S18").Value = 0.1382978723404255
# Row 19
$ws.Range("F19").Value = 0.01979472140762463
$ws.Range("H19").Value = 0.2067448680351906
$ws.Range("I19").Value = 0.08870967741935484
$ws.Range("J19").Value = 0.3724340175953079
$ws.Range("K19").Value = 0.09897360703812316
$ws.Range("M19").Value = 0.02199413489736071
$ws.Range("N19").Value = 0.001466275659824047
$ws.Range("O19").Value = 0.05791788856304985
$ws.Range("S19").Value = 0.1319648093841642
